# Insert a new row above row 47 on the "2024" sheet, shifting the existing
# September SMS-log rows (and everything below them, including the
# "Broadband" label that lived at A197) down by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

$ws.Rows.Item(47).Insert()

# Populate the newly-inserted row 47 with the latest September entry.
$ws.Range("R47").Value = "indusind"
$ws.Range("S47").Value = "2024-09-23 09:08:37"
